$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Mutation Testing Applied to Estelle Specifications
$ws.Range("A7").Value = "Mutation Testing Applied to Estelle Specifications"
$ws.Range("B7").Value = "De Souza S D R S, Maldonado J C, Fabbri S C P F, et al. Mutation testing applied to estelle specifications[J]. Software Quality Journal, 1999, 8(4): 285-301."

# Row 8: Towards practical application of mutation testing in industry (wrap text on A8)
$ws.Range("A8").Value = "Towards practical application of mutation testing in industry — Traditional versus extreme mutation testing"
$ws.Range("A8").WrapText = $true
$ws.Range("B8").Value = "Betka M, Wagner S. Towards practical application of mutation testing in industry—Traditional versus extreme mutation testing[J]. Journal of Software: Evolution and Process, 2022: e2450."

# Row 9: Applying Mutation Testing to XML Schemas
$ws.Range("A9").Value = "Applying Mutation Testing to XML Schemas"
$ws.Range("B9").Value = "Franzotte L, Vergilio S R. Applying Mutation Testing in XML Schemas[C]//SEKE. 2006: 511-516."

# Row 10: Mutation Testing Applied to Hardware: the Mutants Generation
$ws.Range("A10").Value = "Mutation Testing Applied to Hardware: the Mutants Generation"
$ws.Range("B10").Value = "Nguyen T B, Robach C. Mutation testing applied to hardware: the mutants generation[C]//Proceedings of the 11th IFIP International Conference on Very Large Scale Integration. 2001: 118-123."

# Row 11: MDroid+: A Mutation Testing Framework for Android
$ws.Range("A11").Value = "MDroid+: A Mutation Testing Framework for Android"
$ws.Range("B11").Value = "Moran K, Tufano M, Bernal-Cárdenas C, et al. Mdroid+: A mutation testing framework for android[C]//2018 IEEE/ACM 40th International Conference on Software Engineering: Companion (ICSE-Companion). IEEE, 2018: 33-36."

# Update selection to mirror the final workbook state
$ws.Range("B17").Select()
